$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 13: highlight (format-paint from an already-highlighted cell so the
# resulting style exactly matches the existing "yellow" style used elsewhere
# in the sheet, including clearing the pre-existing center alignment on
# B13:D13), then fill in the new Remarks note in H13. ---
$ws.Range("A6").Copy()
$ws.Range("A13:I13").PasteSpecial(-4122)
$ws.Range("H13").Value = "Seen, Not interested"

# --- Row 42: same treatment (I42 is intentionally left untouched). ---
$ws.Range("A6").Copy()
$ws.Range("A42:H42").PasteSpecial(-4122)
$ws.Range("H42").Value = "Seen, Not interested"

# --- Row 47: new listing "Erik Satiestraat 22". ---
$ws.Range("A47").Value = "Erik Satiestraat 22"
$ws.Range("B47").Value = 235
$ws.Range("C47").Value = "1323 SN Almere"
$ws.Range("D47").Value = 114
$ws.Range("E47").Value = "Email Sent"
$ws.Range("F47").Value = "YES"
$ws.Range("G47").Value = "03 July 4PM"
$ws.Range("I47").Value = "https://www.funda.nl/koop/almere/huis-40619868-erik-satiestraat-22/"

# B47:D47 picked up a centered style from the row's default column format;
# paint the plain (non-centered) format from a neighboring cell onto them so
# they match the rest of the row.
$ws.Range("E47").Copy()
$ws.Range("B47:D47").PasteSpecial(-4122)

# --- Restore the current view/selection (scrolled down + H43 selected). ---
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("H43").Select()
